$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.945.89"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.891.78"
$ws.Range("E3").Value = "  -2.42%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.34%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7327"
$ws.Range("E5").Value = "  -2.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.80"
$ws.Range("E6").Value = "  -1.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.37%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3091"
$ws.Range("E8").Value = "  -2.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.20"
$ws.Range("E9").Value = "  -4.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06891"
$ws.Range("E10").Value = "  -1.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7705"
$ws.Range("E11").Value = "  -1.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07949"
$ws.Range("E12").Value = "  -0.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.882.24"
$ws.Range("E13").Value = "  -2.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.213"
$ws.Range("E14").Value = "  -2.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.42"
$ws.Range("E15").Value = "  -3.24%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.959.20"
$ws.Range("E16").Value = "  -1.50%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.10"
$ws.Range("E17").Value = "  -2.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.772"
$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "239.47"
$ws.Range("E19").Value = "  -5.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007756"
$ws.Range("E20").Value = "  -1.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.139.42"
$ws.Range("E22").Value = "  -2.42%  "

$ws.Range("E23").Value = "  -0.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.949"
$ws.Range("E24").Value = "  +4.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.283"
$ws.Range("E25").Value = "  -2.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.25"
$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.81"
$ws.Range("E27").Value = "  -0.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1265"
$ws.Range("E28").Value = "  -4.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.014"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.359"
$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.531"
$ws.Range("E31").Value = "  +1.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.301"
$ws.Range("E32").Value = "  -1.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.056"
$ws.Range("E33").Value = "  -0.99%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05097"
$ws.Range("E34").Value = "  -0.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.275"
$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7339"
$ws.Range("E36").Value = "  -1.57%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.718"
$ws.Range("E37").Value = "  -2.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01920"
$ws.Range("E38").Value = "  -1.25%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.303"
$ws.Range("E40").Value = "  -1.80%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.08"
$ws.Range("E41").Value = "  -5.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4436"
$ws.Range("E42").Value = "  -0.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.929"
$ws.Range("E43").Value = "  -1.68%  "

$ws.Range("E44").Value = "  -0.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8378"
$ws.Range("E45").Value = "  +0.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.617"
$ws.Range("E46").Value = "  +2.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.98"
$ws.Range("E47").Value = "  -0.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.790"
$ws.Range("E48").Value = "  +0.60%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.043.29"
$ws.Range("E49").Value = "  -2.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.33"
$ws.Range("E50").Value = "  -2.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "936.20"
$ws.Range("E51").Value = "  -4.70%  "
